$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns (I1, J1), copying the existing header
# formatting (bold, centered, bordered) from H1 so the new headers
# match the look of the existing ones.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Add the corresponding data values on row 2 (plain, unstyled numbers,
# like the rest of the numeric cells in that row).
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
